$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header B1: "groups" -> "group"
$ws.Range("B1").Value = "group"

# Update row 2: group name and phone number (order matters for shared-string ids)
$ws.Range("B2").Value = "Smiddle"
$ws.Range("A2").Value = "0638845771"

# Delete row 3 (was: 0637926099 / super_VIP)
$ws.Rows("3").Delete()

# Move selection to L7 (matches final selection in the file)
$ws.Range("L7").Select()
